$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("martin power law")
$ws.Activate()

# --- Block 1 (rows 1-10): header row ---
$ws.Range("C1").Value = "Fz C flux organic"
$ws.Range("D1").Value = "fz0(z/z0) Normalized to 102 m"

# Recompute column D (normalized flux) using the true 100 m correction:
# fz0(z/z0) = C$2 * (B_row / B$2)   -- instead of the old C_row / C$2
$ws.Range("D2").Formula  = '=$C$2*(B2/$B$2)'
$ws.Range("D3").Formula  = '=$C$2*(B3/$B$2)'
$ws.Range("D4").Formula  = '=$C$2*(B4/$B$2)'
$ws.Range("D5").Formula  = '=$C$2*(B5/$B$2)'
$ws.Range("D6").Formula  = '=$C$2*(B6/$B$2)'
$ws.Range("D7").Formula  = '=$C$2*(B7/$B$2)'
$ws.Range("D8").Formula  = '=$C$2*(B8/$B$2)'
$ws.Range("D9").Formula  = '=$C$2*(B9/$B$2)'
$ws.Range("D10").Formula = '=$C$2*(B10/$B$2)'

# --- Block 2 (rows 12-24): second header row ---
$ws.Range("C12").Value = "Fz C flux organic"
$ws.Range("D12").Value = "fz0(z/z0) Normalized to 100 m"

# Recompute column D (normalized flux) using the true 100 m correction:
$ws.Range("D13").Formula = '=$C$14*(B13/$B$14)'
$ws.Range("D14").Formula = '=$C$14*(B14/$B$14)'
$ws.Range("D15").Formula = '=$C$14*(B15/$B$14)'
$ws.Range("D16").Formula = '=$C$14*(B16/$B$14)'
$ws.Range("D17").Formula = '=$C$14*(B17/$B$14)'
$ws.Range("D18").Formula = '=$C$14*(B18/$B$14)'
$ws.Range("D19").Formula = '=$C$14*(B19/$B$14)'
$ws.Range("D20").Formula = '=$C$14*(B20/$B$14)'
$ws.Range("D21").Formula = '=$C$14*(B21/$B$14)'
$ws.Range("D22").Formula = '=$C$14*(B22/$B$14)'
$ws.Range("D23").Formula = '=$C$14*(B23/$B$14)'
$ws.Range("D24").Formula = '=$C$14*(B24/$B$14)'

# Restore the cursor / scroll position recorded in the sheet view.
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("S39").Select() | Out-Null
